$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EoIEECwEC")

$ws.Range("A12").Value = "green hydrogen if"
$ws.Range("A13").Value = "low carbon hydrogen if"

$ws.Range("B12:B13").Formula = "=B3"

$ws.Activate()
$ws.Range("A14").Select() | Out-Null
